$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 0.074444
$ws.Range("M2").Value = 19.77408333333333
$ws.Range("N2").Value = 59.32225
$ws.Range("O2").Value = 0.3380388258879848
$ws.Range("P2").Value = 0.339186328349942
$ws.Range("Q2").Value = 0.4906872865555555
$ws.Range("R2").Value = 4.416185579
$ws.Range("S2").Value = 0.3380388258879848
$ws.Range("T2").Value = 0.339186328349942

# Row 3
$ws.Range("H3").Value = 0.074444
$ws.Range("O3").Value = 0.3069959581674471
$ws.Range("P3").Value = 0.3080380828904952
$ws.Range("R3").Value = 4.010637297976
$ws.Range("S3").Value = 0.3069959581674471
$ws.Range("T3").Value = 0.3080380828904952

# Row 4
$ws.Range("H4").Value = 0.074444
$ws.Range("M4").Value = 12.46730333333333
$ws.Range("N4").Value = 37.40191
$ws.Range("O4").Value = 0.2131290998296268
$ws.Range("P4").Value = 0.2138525852639604
$ws.Range("Q4").Value = 0.3093719764488889
$ws.Range("R4").Value = 2.78434778804
$ws.Range("S4").Value = 0.2131290998296268
$ws.Range("T4").Value = 0.2138525852639604

# Row 5
$ws.Range("H5").Value = 0.074444
$ws.Range("M5").Value = 0.593699
$ws.Range("N5").Value = 1.187398
$ws.Range("O5").Value = 0.01014931056513554
$ws.Range("P5").Value = 0.006789175527058808
$ws.Range("Q5").Value = 0.01473244278533333
$ws.Range("R5").Value = 0.08839465671199999
$ws.Range("S5").Value = 0.01014931056513554
$ws.Range("T5").Value = 0.006789175527058808

# Row 6
$ws.Range("H6").Value = 0.074444
$ws.Range("M6").Value = 7.703215333333333
$ws.Range("N6").Value = 23.109646
$ws.Range("O6").Value = 0.1316868055498057
$ws.Range("P6").Value = 0.1321338279685434
$ws.Range("Q6").Value = 0.1911527207582222
$ws.Range("R6").Value = 1.720374486824
$ws.Range("S6").Value = 0.1316868055498057
$ws.Range("T6").Value = 0.1321338279685434
